# Regenerate merged AHB files
# - Rename the "_old" / "_new" header-suffix columns to the version-specific
#   "_FV2410" / "_FV2504" suffixes.
# - Freeze the header row.
# - Turn the A1:U73 range into an Excel Table (ListObject).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename header row (row 1) -----------------------------------------
# Columns A:J carried the "_old" suffix -> "_FV2410"
# Column K is the unchanged "diff" column
# Columns L:U carried the "_new" suffix -> "_FV2504"
$headers = @(
    "Segmentname_FV2410",
    "Segmentgruppe_FV2410",
    "Segment_FV2410",
    "Datenelement_FV2410",
    "Segment ID_FV2410",
    "Code_FV2410",
    "Qualifier_FV2410",
    "Beschreibung_FV2410",
    "Bedingungsausdruck_FV2410",
    "Bedingung_FV2410",
    "diff",
    "Segmentname_FV2504",
    "Segmentgruppe_FV2504",
    "Segment_FV2504",
    "Datenelement_FV2504",
    "Segment ID_FV2504",
    "Code_FV2504",
    "Qualifier_FV2504",
    "Beschreibung_FV2504",
    "Bedingungsausdruck_FV2504",
    "Bedingung_FV2504"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# --- 2. Freeze the header row -----------------------------------------
$ws.Range("A2").Select() | Out-Null
[void]($excel.ActiveWindow.FreezePanes = $true)

# --- 3. Convert the used range into a table ----------------------------
$tbl = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $ws.Range("A1:U73"), $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$tbl.Name = "Table1"
